# close #187: Remove unnecessary column name in values and proportionality
# The worksheet had a "nome" (name) column in column B that is no longer
# needed; delete the entire column B and shift everything left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()

# Leave the selection on B1, matching the post-edit state.
$ws.Range("B1").Select()
